# OLX Monitor 2026-02-23 09:40 — append the newest scrape results to the
# "PODSUMOWANIE" sheet's running log (rows 211-218), mirroring the exact
# per-column style pattern already used by the rest of the log.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Clone formatting from the last existing data row (210) onto the
#        8 new rows, so they pick up the same per-column styles (left/
#        center alignment, the "stale listing" red-font style in F, etc.)
#        that the rest of the log already uses.
$srcRow = $ws.Range("A210:H210")
$dstRows = $ws.Range("A211:H218")
$srcRow.Copy()
$dstRows.PasteSpecial(-4122)

# --- 2. Two of the new rows ("fresh" listings re-seen at 34 days) use the
#        plain (non-red) style in column F instead of the "stale" red
#        style copied above. Row 204 already carries that exact style, so
#        borrow it for F212 / F216.
$fNormal = $ws.Range("F204")
$fNormal.Copy()
$ws.Range("F212").PasteSpecial(-4122)
$ws.Range("F216").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Row data (Profile / timestamp / title / price / date / days /
#        url / slug) for the 8 listings seen in this run.
$rows = @(
    @{ r = 211; a = "2026-02-23 09:40:00"; b = "poqui";           c = "Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza";                            d = 2049;  e = "19.12.2025"; f = 65;  g = "https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html";                       h = "mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc" }
    @{ r = 212; a = "2026-02-23 09:40:00"; b = "poqui";           c = "Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda";                                 d = 2299;  e = "19.01.2026"; f = 34;  g = "https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html";                             h = "swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR" }
    @{ r = 213; a = "2026-02-23 09:40:00"; b = "poqui";           c = "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy";                            d = 2499;  e = "28.10.2025"; f = 117; g = "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html";                     h = "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger" }
    @{ r = 214; a = "2026-02-23 09:40:00"; b = "poqui";           c = "Przytulny pokój blisko Politechniki – ul. Przytulna";                                      d = 549;   e = "10.10.2025"; f = 136; g = "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html";                             h = "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz" }
    @{ r = 215; a = "2026-02-23 09:40:00"; b = "pokojewlublinie"; c = "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58";                                      d = 0;     e = "11.08.2025"; f = 195; g = "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html";                             h = "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm" }
    @{ r = 216; a = "2026-02-23 09:40:00"; b = "pokojewlublinie"; c = "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12";                    d = 12640; e = "19.01.2026"; f = 34;  g = "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html";       h = "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc" }
    @{ r = 217; a = "2026-02-23 09:40:00"; b = "dawnypatron";     c = "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.";                    d = 730;   e = "20.09.2024"; f = 520; g = "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html";             h = "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM" }
    @{ r = 218; a = "2026-02-23 09:40:00"; b = "dawnypatron";     c = "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14";                    d = 14690; e = "05.12.2025"; f = 79;  g = "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html"; h = "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv" }
)

# Day-month-ambiguous dates (dd <= 12) get auto-parsed into real Excel
# dates by the smart Value setter. Force those particular cells to text
# first so the literal "dd.mm.yyyy" string survives.
$ambiguousDateRows = @(214, 215, 218)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.a
    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $row.d

    if ($ambiguousDateRows -contains $r) {
        $ws.Cells.Item($r, 5).NumberFormat = "@"
        $ws.Cells.Item($r, 5).Value = $row.e
    } else {
        $ws.Cells.Item($r, 5).Value = $row.e
    }

    $ws.Cells.Item($r, 6).Value = $row.f
    $ws.Cells.Item($r, 7).Value = $row.g
    $ws.Cells.Item($r, 8).Value = $row.h
}

# Re-paste the column-E format over the forced-text cells so the
# NumberFormat goes back to General (matching the rest of the log) while
# the already-stored text value is left untouched.
$eNormal = $ws.Range("E210")
$eNormal.Copy()
foreach ($r in $ambiguousDateRows) {
    $ws.Cells.Item($r, 5).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

Write-Output "Appended rows 211-218 to $($ws.Name)"
